# Apply the "reorder git metadata attributes (url, branch, revision)" edit
# to the "Data repo metadata" sheet, add a trailing blank row, relocate the
# data-validation ranges to match, and move the active-tab/selection from
# "Model1s" to "Data repo metadata".

$wb = $excel.ActiveWorkbook

$wsData  = $wb.Worksheets.Item("Data repo metadata")
$wsModel = $wb.Worksheets.Item("Model1s")

# --- 1. Re-order the metadata rows on "Data repo metadata" ---------------
# Before: row1=Branch/master, row2=Revision/<sha>, row3=Url/<url>
# After:  row1=Url/<url>,     row2=Branch/master,  row3=Revision/<sha>
$branchLabel = $wsData.Range("A1").Value2
$branchValue = $wsData.Range("B1").Value2
$revisionLabel = $wsData.Range("A2").Value2
$revisionValue = $wsData.Range("B2").Value2
$urlLabel = $wsData.Range("A3").Value2
$urlValue = $wsData.Range("B3").Value2

$wsData.Range("A1").Value = $urlLabel
$wsData.Range("B1").Value = $urlValue
$wsData.Range("A2").Value = $branchLabel
$wsData.Range("B2").Value = $branchValue
$wsData.Range("A3").Value = $revisionLabel
$wsData.Range("B3").Value = $revisionValue

# --- 2. Move the data validations so they still point at the right label -
$wsData.Range("B1").Validation.Delete()
$wsData.Range("B2").Validation.Delete()
$wsData.Range("B3").Validation.Delete()

# Url validation -> B1
$wsData.Range("B1").Validation.Add(6, 1, 8, 255)
$wsData.Range("B1").Validation.ErrorTitle = "Url"
$wsData.Range("B1").Validation.ErrorMessage = "Value must be a string._x000a__x000a_Value must be less than or equal to 255 characters."
$wsData.Range("B1").Validation.InputTitle = "Url"
$wsData.Range("B1").Validation.InputMessage = "Enter a string._x000a__x000a_Value must be less than or equal to 255 characters."
$wsData.Range("B1").Validation.IgnoreBlank = $true
$wsData.Range("B1").Validation.InCellDropdown = $true
$wsData.Range("B1").Validation.ShowInput = $true
$wsData.Range("B1").Validation.ShowError = $true
$wsData.Range("B1").Validation.AlertStyle = 2

# Branch validation -> B2
$wsData.Range("B2").Validation.Add(6, 1, 8, 255)
$wsData.Range("B2").Validation.ErrorTitle = "Branch"
$wsData.Range("B2").Validation.ErrorMessage = "Value must be a string._x000a__x000a_Value must be less than or equal to 255 characters."
$wsData.Range("B2").Validation.InputTitle = "Branch"
$wsData.Range("B2").Validation.InputMessage = "Enter a string._x000a__x000a_Value must be less than or equal to 255 characters."
$wsData.Range("B2").Validation.IgnoreBlank = $true
$wsData.Range("B2").Validation.InCellDropdown = $true
$wsData.Range("B2").Validation.ShowInput = $true
$wsData.Range("B2").Validation.ShowError = $true
$wsData.Range("B2").Validation.AlertStyle = 2

# Revision validation -> B3
$wsData.Range("B3").Validation.Add(6, 1, 8, 255)
$wsData.Range("B3").Validation.ErrorTitle = "Revision"
$wsData.Range("B3").Validation.ErrorMessage = "Value must be a string._x000a__x000a_Value must be less than or equal to 255 characters."
$wsData.Range("B3").Validation.InputTitle = "Revision"
$wsData.Range("B3").Validation.InputMessage = "Enter a string._x000a__x000a_Value must be less than or equal to 255 characters."
$wsData.Range("B3").Validation.IgnoreBlank = $true
$wsData.Range("B3").Validation.InCellDropdown = $true
$wsData.Range("B3").Validation.ShowInput = $true
$wsData.Range("B3").Validation.ShowError = $true
$wsData.Range("B3").Validation.AlertStyle = 2

# --- 3. Add the new (blank) trailing row 4 --------------------------------
$wsData.Range("A4").Font.Bold = $false
$wsData.Rows.Item(4).RowHeight = 15

# --- 4. "Model1s" loses tabSelected and gets a new zoom level -------------
$wsModel.Activate()
$win = $wb.Windows.Item(1)
$win.Zoom = 130

# --- 5. Move the active tab / selection to "Data repo metadata"
# (activeTab goes from 3 -> 1, i.e. "Model1s" -> "Data repo metadata").
$wsData.Activate()
$wsData.Range("A4:XFD4").Select()
